# Fruta / hortaliza, semanal
# Insert a new weekly block of 3 rows (Especial/Primera/Segunda) for Mango
# at the top of the existing date-ordered data (rows 1177-1179), pushing
# all subsequent rows down by 3 (one block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current row 1177, shifting 1177:1224 -> 1180:1227
$ws.Rows.Item(1177).Insert()
$ws.Rows.Item(1178).Insert()
$ws.Rows.Item(1179).Insert()

# Common (static) values shared by the whole new block, copied from the
# surrounding Mango / Terminal La Palmera de La Serena rows.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$prodId    = 100108
$producto  = "Tropicales y subtropicales"
$catId     = 100108002
$categoria = "Mango"
$variedad  = "Sin especificar"
$volumen   = 512
$unidad    = "`$/bandeja 4 kilos"
$kgUnidad  = 4

$fecha   = 45075
$pmin    = 7500
$pmax    = 8000
$pprom   = 7750
$origen  = "Perú"
$pkg     = 1938

$calidades = @("Especial", "Primera", "Segunda")

for ($i = 0; $i -lt 3; $i++) {
    $r = 1177 + $i
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $prodId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidades[$i]
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $pkg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
